$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Remove the old extra rows (4-8); rows 1-3 will be overwritten below ---
$ws.Rows("4:8").Delete() | Out-Null

# --- Header row ---
$ws.Cells.Item(1,1).Value = "Date"
$ws.Cells.Item(1,2).Value = "Rohit Malvi"
$ws.Cells.Item(1,3).Value = "Amol Aldar"
$ws.Cells.Item(1,4).Value = "Shubham Sanap"
$ws.Cells.Item(1,5).Value = "Pravin Mane"
$ws.Cells.Item(1,6).Value = "Prathmesh Patil"
$ws.Cells.Item(1,7).Value = "Kiran Kale"
$ws.Cells.Item(1,8).Value = "Amar Maurya"
$ws.Cells.Item(1,9).Value = "Raghupati Khot"

# --- Row 2 (1st Aug 2023) ---
$ws.Cells.Item(2,1).Value = 45139
$ws.Cells.Item(2,1).NumberFormat = "d-mmm-yy"
$ws.Cells.Item(2,2).Value = "Present"
$ws.Cells.Item(2,3).Value = "Present"
$ws.Cells.Item(2,4).Value = "Present"
$ws.Cells.Item(2,5).Value = "Present"
$ws.Cells.Item(2,6).Value = "Absent"
$ws.Cells.Item(2,7).Value = "Present"
$ws.Cells.Item(2,8).Value = "Present"
$ws.Cells.Item(2,9).Value = "Absent"

# --- Row 3 (2nd Aug 2023) ---
$ws.Cells.Item(3,1).Value = 45140
$ws.Cells.Item(3,1).NumberFormat = "d-mmm-yy"
$ws.Cells.Item(3,2).Value = "Present"
$ws.Cells.Item(3,3).Value = "Present"
$ws.Cells.Item(3,4).Value = "Present"
$ws.Cells.Item(3,5).Value = "Absent"
$ws.Cells.Item(3,6).Value = "Present"
$ws.Cells.Item(3,7).Value = "Present"
$ws.Cells.Item(3,8).Value = "Present"
$ws.Cells.Item(3,9).Value = "Present"

# --- Column widths (engine quantizes ColumnWidth to 1/6-character steps, so
#     these inputs are chosen to land as close as possible on the target
#     stored OOXML widths of 15, 14.42578125, 18.42578125, 12.7109375, 15.85546875) ---
$ws.Columns("D").ColumnWidth = 14.166666666666666
$ws.Columns("E").ColumnWidth = 13.592447916666666
$ws.Columns("G").ColumnWidth = 17.592447916666668
$ws.Columns("H").ColumnWidth = 11.877604166666666
$ws.Columns("I").ColumnWidth = 15.022135416666666

# --- Comments ---
$c1 = $ws.Range("I2").AddComment("HP:`nHe was unable to join call as he was in office.")
$c2 = $ws.Range("E3").AddComment("HP:`nDue to personal reason he was unable to join the session.")

# --- Selection ---
$ws.Range("A4").Select() | Out-Null

Write-Output "done"
